# Updates the "cryptos" worksheet with refreshed Price (column D) and
# Volume(1h) (column E) figures, mirroring the GitHub Actions data refresh
# captured in the commit "Updated cryptos list ... with GitHub Actions".
#
# Only rows 2-51 are touched; column D is left untouched where the source
# diff shows no change to the price text (rows 10, 28, 30, 33, 41, 42, 48, 49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '66.505.17'; E = '  -0.84%  ' },
    @{ Row = 3; D = '3.079.63'; E = '  -0.98%  ' },
    @{ Row = 4; D = '0.999'; E = '  -0.07%  ' },
    @{ Row = 5; D = '574.43'; E = '  -0.75%  ' },
    @{ Row = 6; D = '176.59'; E = '  +1.86%  ' },
    @{ Row = 7; D = '0.999'; E = '  -0.05%  ' },
    @{ Row = 8; D = '3.079.16'; E = '  -0.89%  ' },
    @{ Row = 9; D = '0.511'; E = '  -1.82%  ' },
    @{ Row = 10; D = $null; E = '  -2.77%  ' },
    @{ Row = 11; D = '0.150'; E = '  -2.54%  ' },
    @{ Row = 12; D = '0.463'; E = '  -3.11%  ' },
    @{ Row = 13; D = '0.0000238'; E = '  -3.69%  ' },
    @{ Row = 14; D = '35.72'; E = '  -3.13%  ' },
    @{ Row = 15; D = '0.121'; E = '  -0.63%  ' },
    @{ Row = 16; D = '3.592.37'; E = '  -0.85%  ' },
    @{ Row = 17; D = '66.415.05'; E = '  -0.91%  ' },
    @{ Row = 18; D = '6.93'; E = '  -2.18%  ' },
    @{ Row = 19; D = '16.75'; E = '  +1.93%  ' },
    @{ Row = 20; D = '3.076.88'; E = '  -1.03%  ' },
    @{ Row = 21; D = '478.23'; E = '  -2.83%  ' },
    @{ Row = 22; D = '7.69'; E = '  -2.85%  ' },
    @{ Row = 23; D = '0.684'; E = '  -2.81%  ' },
    @{ Row = 24; D = '83.08'; E = '  -0.96%  ' },
    @{ Row = 25; D = '12.57'; E = '  -4.31%  ' },
    @{ Row = 26; D = '2.20'; E = '  -3.79%  ' },
    @{ Row = 27; D = '10.03'; E = '  -4.14%  ' },
    @{ Row = 28; D = $null; E = '  +0.08%  ' },
    @{ Row = 29; D = '7.95'; E = '  +0.62%  ' },
    @{ Row = 30; D = $null; E = '  -4.33%  ' },
    @{ Row = 31; D = '2.58'; E = '  -3.07%  ' },
    @{ Row = 32; D = '27.81'; E = '  -2.07%  ' },
    @{ Row = 33; D = $null; E = '  -2.36%  ' },
    @{ Row = 34; D = "0.0$([char]0x2083)0935"; E = '  -1.29%  ' },
    @{ Row = 35; D = '0.998'; E = '  -0.10%  ' },
    @{ Row = 36; D = '47.91'; E = '  +1.62%  ' },
    @{ Row = 37; D = '5.54'; E = '  -5.36%  ' },
    @{ Row = 38; D = '0.936'; E = '  -3.72%  ' },
    @{ Row = 39; D = '48.88'; E = '  -2.29%  ' },
    @{ Row = 40; D = '0.307'; E = '  -0.43%  ' },
    @{ Row = 41; D = $null; E = '  -0.79%  ' },
    @{ Row = 42; D = $null; E = '  -3.57%  ' },
    @{ Row = 43; D = '8.27'; E = '  -2.22%  ' },
    @{ Row = 44; D = '2.65'; E = '  +2.29%  ' },
    @{ Row = 45; D = '2.775.28'; E = '  -1.03%  ' },
    @{ Row = 46; D = '369.59'; E = '  -3.95%  ' },
    @{ Row = 47; D = '135.38'; E = '  -0.18%  ' },
    @{ Row = 48; D = $null; E = '  -2.58%  ' },
    @{ Row = 49; D = $null; E = '  +0.00%  ' },
    @{ Row = 50; D = '24.76'; E = '  -0.41%  ' },
    @{ Row = 51; D = '2.20'; E = '  +0.54%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
